$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RVL")

# The "Param Value" column (G) for the two "Global / DoSleep / millis / number"
# rows (rows 3 and 5) changes from 2000 to 1000. The value must stay a text
# cell (it was authored as shared-string text, not a number), so it's entered
# with a leading apostrophe to force text. That alone would also flip Excel's
# internal "quote prefix" formatting flag on the cell, so we immediately
# paste-special just the number format back from the untouched neighboring
# cell (column F, same row) to restore the original (default) cell style
# while keeping the new text value.
$xlPasteFormats = -4122

$ws.Range("G3").Value = "'1000"
$ws.Range("F3").Copy()
$ws.Range("G3").PasteSpecial($xlPasteFormats)

$ws.Range("G5").Value = "'1000"
$ws.Range("F5").Copy()
$ws.Range("G5").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false
